$d = $word.ActiveDocument

$replacements = @(
    @("2025-09-13 Saturday", "2025-09-14 Sunday"),
    @("57×14=798", "84×43=3612"),
    @("34×37=1258", "75×77=5775"),
    @("26×33=858", "69×86=5934"),
    @("51×88=4488", "74×34=2516"),
    @("63×71=4473", "40×27=1080"),
    @("11×68=748", "64×41=2624"),
    @("31×21=651", "37×40=1480"),
    @("73×24=1752", "16×28=448"),
    @("33×66=2178", "41×85=3485"),
    @("92×74=6808", "48×63=3024"),
    @("35×57=1995", "23×19=437"),
    @("43×37=1591", "79×86=6794"),
    @("79×64=5056", "13×20=260"),
    @("87×46=4002", "11×66=726"),
    @("20×90=1800", "57×12=684"),
    @("35×96=3360", "61×69=4209"),
    @("73×13=949", "51×13=663"),
    @("92×50=4600", "19×30=570"),
    @("38×29=1102", "45×17=765"),
    @("19×76=1444", "38×78=2964"),
    @("49×71=3479", "59×66=3894"),
    @("20×45=900", "21×71=1491"),
    @("27×12=324", "21×87=1827"),
    @("49×91=4459", "11×18=198"),
    @("84×27=2268", "54×50=2700")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
